# Scheduled-runner update: refresh computed leve-profit figures (H/I/J/K/L/M/N
# columns) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with the latest
# market-board snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3416.9412
$ws.Cells.Item(40, 9).Value = 3525
$ws.Cells.Item(40, 11).Value = 3525
$ws.Cells.Item(40, 13).Value = -3350
$ws.Cells.Item(42, 8).Value = 714.25
$ws.Cells.Item(42, 9).Value = 319
$ws.Cells.Item(42, 11).Value = 957
$ws.Cells.Item(42, 13).Value = -727
$ws.Cells.Item(62, 8).Value = 9954.689
$ws.Cells.Item(62, 9).Value = 3783
$ws.Cells.Item(62, 10).Value = 16567.215
$ws.Cells.Item(62, 11).Value = 3783
$ws.Cells.Item(62, 12).Value = 16567.215
$ws.Cells.Item(62, 13).Value = -3159
$ws.Cells.Item(62, 14).Value = -17815.215
$ws.Cells.Item(65, 8).Value = 9954.689
$ws.Cells.Item(65, 9).Value = 3783
$ws.Cells.Item(65, 10).Value = 16567.215
$ws.Cells.Item(65, 11).Value = 18915
$ws.Cells.Item(65, 12).Value = 82836.075
$ws.Cells.Item(65, 13).Value = -15795
$ws.Cells.Item(65, 14).Value = -89076.075
$ws.Cells.Item(105, 8).Value = 39245
$ws.Cells.Item(105, 10).Value = 39245
$ws.Cells.Item(105, 12).Value = 39245
$ws.Cells.Item(105, 14).Value = -46233
$ws.Cells.Item(112, 8).Value = 2148.6667
$ws.Cells.Item(112, 9).Value = 1773.3334
$ws.Cells.Item(112, 11).Value = 5320.0002
$ws.Cells.Item(112, 13).Value = -4212.0002
$ws.Cells.Item(132, 8).Value = 1972.8064
$ws.Cells.Item(132, 9).Value = 1048.1923
$ws.Cells.Item(132, 11).Value = 3144.5769
$ws.Cells.Item(132, 13).Value = -614.5769

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8083309.5
$ws.Cells.Item(32, 9).Value = 9109585
$ws.Cells.Item(32, 10).Value = 19716.285
$ws.Cells.Item(32, 11).Value = 9109585
$ws.Cells.Item(32, 12).Value = 19716.285
$ws.Cells.Item(32, 13).Value = -9109298
$ws.Cells.Item(32, 14).Value = -20290.285
$ws.Cells.Item(39, 8).Value = 3681
$ws.Cells.Item(39, 9).Value = 3681
$ws.Cells.Item(39, 11).Value = 3681
$ws.Cells.Item(39, 13).Value = -3161
$ws.Cells.Item(45, 8).Value = 27779872
$ws.Cells.Item(45, 9).Value = 35716092
$ws.Cells.Item(45, 10).Value = 3101
$ws.Cells.Item(45, 11).Value = 35716092
$ws.Cells.Item(45, 12).Value = 3101
$ws.Cells.Item(45, 13).Value = -35715715
$ws.Cells.Item(45, 14).Value = -3855
$ws.Cells.Item(61, 8).Value = 107154290
$ws.Cells.Item(61, 9).Value = 250010000
$ws.Cells.Item(61, 11).Value = 250010000
$ws.Cells.Item(61, 13).Value = -250009788
$ws.Cells.Item(123, 8).Value = 72000
$ws.Cells.Item(123, 10).Value = 72000
$ws.Cells.Item(123, 12).Value = 72000
$ws.Cells.Item(123, 14).Value = -81800
$ws.Cells.Item(136, 8).Value = 107154290
$ws.Cells.Item(136, 9).Value = 250010000
$ws.Cells.Item(136, 11).Value = 750030000
$ws.Cells.Item(136, 13).Value = -750027450

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 16001567
$ws.Cells.Item(7, 9).Value = 1959.25
$ws.Cells.Item(7, 10).Value = 80000000
$ws.Cells.Item(7, 11).Value = 1959.25
$ws.Cells.Item(7, 12).Value = 80000000
$ws.Cells.Item(7, 13).Value = -1846.25
$ws.Cells.Item(7, 14).Value = -80000226
$ws.Cells.Item(26, 8).Value = 29823.666
$ws.Cells.Item(26, 9).Value = 19735.5
$ws.Cells.Item(26, 11).Value = 19735.5
$ws.Cells.Item(26, 13).Value = -19443.5
$ws.Cells.Item(38, 8).Value = 45047.07
$ws.Cells.Item(38, 10).Value = 38291.8
$ws.Cells.Item(38, 12).Value = 38291.8
$ws.Cells.Item(38, 14).Value = -39123.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 3868.8572
$ws.Cells.Item(107, 9).Value = 3830.3333
$ws.Cells.Item(107, 10).Value = 4100
$ws.Cells.Item(107, 11).Value = 3830.3333
$ws.Cells.Item(107, 12).Value = 4100
$ws.Cells.Item(107, 13).Value = -1910.3333
$ws.Cells.Item(107, 14).Value = -7940

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(118, 8).Value = 5348.1665
$ws.Cells.Item(118, 10).Value = 8966.666999999999
$ws.Cells.Item(118, 12).Value = 26900.001
$ws.Cells.Item(118, 14).Value = -29386.001
$ws.Cells.Item(141, 8).Value = 197291.12
$ws.Cells.Item(141, 9).Value = 379719.75
$ws.Cells.Item(141, 11).Value = 1139159.25
$ws.Cells.Item(141, 13).Value = -1133979.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 25500
$ws.Cells.Item(70, 8).Value = 5770.7144
$ws.Cells.Item(70, 9).Value = 5098.75
$ws.Cells.Item(70, 11).Value = 5098.75
$ws.Cells.Item(70, 13).Value = -4828.75
$ws.Cells.Item(73, 8).Value = 5770.7144
$ws.Cells.Item(73, 9).Value = 5098.75
$ws.Cells.Item(73, 11).Value = 5098.75
$ws.Cells.Item(73, 13).Value = -4162.75
$ws.Cells.Item(81, 8).Value = 25500
$ws.Cells.Item(84, 8).Value = 25500
$ws.Cells.Item(122, 8).Value = 2468
$ws.Cells.Item(122, 9).Value = 2669
$ws.Cells.Item(122, 11).Value = 8007
$ws.Cells.Item(122, 13).Value = -5557

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1382.4286
$ws.Cells.Item(16, 9).Value = 1382.4286
$ws.Cells.Item(16, 11).Value = 1382.4286
$ws.Cells.Item(16, 13).Value = -1212.4286
$ws.Cells.Item(22, 8).Value = 1825.9524
$ws.Cells.Item(22, 9).Value = 1878.9166
$ws.Cells.Item(22, 10).Value = 1755.3334
$ws.Cells.Item(22, 11).Value = 1878.9166
$ws.Cells.Item(22, 12).Value = 1755.3334
$ws.Cells.Item(22, 13).Value = -1583.9166
$ws.Cells.Item(22, 14).Value = -2345.3334
$ws.Cells.Item(27, 8).Value = 1825.9524
$ws.Cells.Item(27, 9).Value = 1878.9166
$ws.Cells.Item(27, 10).Value = 1755.3334
$ws.Cells.Item(27, 11).Value = 1878.9166
$ws.Cells.Item(27, 12).Value = 1755.3334
$ws.Cells.Item(27, 13).Value = -1771.9166
$ws.Cells.Item(27, 14).Value = -1969.3334
$ws.Cells.Item(40, 8).Value = 5204.2856
$ws.Cells.Item(40, 9).Value = 4844.727
$ws.Cells.Item(40, 10).Value = 5599.8
$ws.Cells.Item(40, 11).Value = 4844.727
$ws.Cells.Item(40, 12).Value = 5599.8
$ws.Cells.Item(40, 13).Value = -4708.727
$ws.Cells.Item(40, 14).Value = -5871.8
$ws.Cells.Item(46, 8).Value = 2742.7856
$ws.Cells.Item(46, 9).Value = 1914.5714
$ws.Cells.Item(46, 10).Value = 3571
$ws.Cells.Item(46, 11).Value = 1914.5714
$ws.Cells.Item(46, 12).Value = 3571
$ws.Cells.Item(46, 13).Value = -1726.5714
$ws.Cells.Item(46, 14).Value = -3947
$ws.Cells.Item(93, 8).Value = 29413224
$ws.Cells.Item(93, 9).Value = 50001216
$ws.Cells.Item(93, 10).Value = 1805.7858
$ws.Cells.Item(93, 11).Value = 50001216
$ws.Cells.Item(93, 12).Value = 1805.7858
$ws.Cells.Item(93, 13).Value = -49999968
$ws.Cells.Item(93, 14).Value = -4301.7858

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 395000
$ws.Cells.Item(5, 9).Value = 750000
$ws.Cells.Item(5, 11).Value = 750000
$ws.Cells.Item(5, 13).Value = -749888
$ws.Cells.Item(59, 8).Value = 45000
$ws.Cells.Item(59, 10).Value = 45000
$ws.Cells.Item(59, 12).Value = 45000
$ws.Cells.Item(59, 14).Value = -46476
$ws.Cells.Item(107, 8).Value = 20834408
$ws.Cells.Item(107, 10).Value = 1272.8334
$ws.Cells.Item(107, 12).Value = 3818.5002
$ws.Cells.Item(107, 14).Value = -7658.5002
$ws.Cells.Item(124, 8).Value = 96492.25
$ws.Cells.Item(124, 10).Value = 96492.25
$ws.Cells.Item(124, 12).Value = 96492.25
$ws.Cells.Item(124, 14).Value = -106312.25
$ws.Cells.Item(141, 8).Value = 61999.668
$ws.Cells.Item(141, 10).Value = 61999.668
$ws.Cells.Item(141, 12).Value = 61999.668
$ws.Cells.Item(141, 14).Value = -72359.66800000001
